$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is a maze/map grid. Rows 3-31, columns C:AD ("path" cells that
# alternated between "wall" (value 0) and "path" (value 3)) are all filled
# in to become "path" cells: value 3 with the green highlight fill used
# elsewhere in the sheet (matches existing style s="4" -> fill FF92D050).
$range = $ws.Range("C3:AD31")
$range.Value = 3
$range.Interior.Color = 5296274   # RGB(146, 208, 80) == 0xFF92D050, packed as BGR for OLE COLORREF

# Update the saved selection/view: the author scrolled back up (no frozen
# topLeftCell override) and left the active cell on AV19.
[void]$ws.Range("AV19").Select()
